# Apply "Linking, shapes and key" edits to the Book of Magnus reading-order
# workbook: fill in the "Follows" column (D) for a number of titles and
# normalize the cell style of the ones that already were / now are filled in,
# then update the selection to reflect where the author ended up working.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Fill in / correct the "Follows" values -------------------------------
# (processed in an order that reproduces the exact shared-string table
#  ordering seen in the target workbook)
$ws.Range("D19").Value = "Legion,Fulgrim"
$ws.Range("D22").Value = "The First Heretic,Battle for the Abyss"
$ws.Range("D28").Value = "Know No Fear,Betrayer"
$ws.Range("D31").Value = "Battle for the Abyss,Vulkan Lives,Fallen Angels"
$ws.Range("D29").Value = "Fulgrim,Legion"
$ws.Range("D30").Value = "Deliverance Lost"
$ws.Range("D25").Value = "Fulgrim"
$ws.Range("D33").Value = "Fulgrim"
$ws.Range("D36").Value = "Corax : Soulforge"
$ws.Range("D42").Value = "The Unremembered Empire"
$ws.Range("D47").Value = "The Unremembered Empire"
$ws.Range("D51").Value = "Pharos"

# --- Normalize cell formatting ---------------------------------------------
# Most "Follows" cells use the plain (non-centered) text style already used
# elsewhere in column D; copy it across, one cell at a time, so each of them
# picks up that style (pasting format onto a multi-area range only affects
# the first area, so loop instead of using a single union Range).
$plainCells = @("D9","D19","D22","D25","D28","D29","D30","D31","D33","D42","D47","D51")
foreach ($cellRef in $plainCells) {
    $ws.Range("D4").Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

# D36 ("Ravenlord" follows the novella "Corax : Soulforge") gets the blue
# "novella" style used elsewhere on the sheet instead.
$ws.Range("A20").Copy()
$ws.Range("D36").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Update the view/selection to match where editing finished -------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 16
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("D29").Select()
